$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1, 9).Value = 'Other found locations'

$ws.Cells.Item(2, 5).Value = '[Ying%Liang%NULL%1,  Jing-Jin%Liang%NULL%1,  Qing-Tao%Zhou%NULL%1,  Xiao-Guang%Li%NULL%1,  Fei%Lin%NULL%1,  Zhong-Hua%Deng%NULL%1,  Bi-Ying%Zhang%NULL%1,  Lu%Li%NULL%1,  Xiao-Hua%Wang%NULL%1,  Hong%Zhu%NULL%1,  Qing-Bian%Ma%NULL%1,  Xiao-Mei%Tong%NULL%1,  Jie%Xu%NULL%1,  Yong-Chang%Sun%NULL%1,  Pei-Fang%Wei%NULL%2,  Pei-Fang%Wei%NULL%0]'
$ws.Cells.Item(2, 9).Value = '_PMC'

$ws.Cells.Item(3, 5).Value = '[Xiao-Wei%Xu%NULL%0,  Xiao-Xin%Wu%NULL%1,  Xian-Gao%Jiang%NULL%1,  Kai-Jin%Xu%NULL%1,  Ling-Jun%Ying%NULL%1,  Chun-Lian%Ma%NULL%1,  Shi-Bo%Li%NULL%1,  Hua-Ying%Wang%NULL%1,  Sheng%Zhang%NULL%1,  Hai-Nv%Gao%NULL%1,  Ji-Fang%Sheng%NULL%1,  Hong-Liu%Cai%NULL%1,  Yun-Qing%Qiu%NULL%1,  Lan-Juan%Li%NULL%1]'
$ws.Cells.Item(3, 9).Value = '_PMC'

$ws.Cells.Item(4, 5).Value = '[Pavan K.%Bhatraju%NULL%0,  Bijan J.%Ghassemieh%NULL%1,  Michelle%Nichols%NULL%1,  Richard%Kim%NULL%1,  Keith R.%Jerome%NULL%1,  Arun K.%Nalla%NULL%1,  Alexander L.%Greninger%NULL%1,  Sudhakar%Pipavath%NULL%1,  Mark M.%Wurfel%NULL%1,  Laura%Evans%NULL%1,  Patricia A.%Kritek%NULL%1,  T. Eoin%West%NULL%1,  Andrew%Luks%NULL%1,  Anthony%Gerbino%NULL%1,  Chris R.%Dale%NULL%1,  Jason D.%Goldman%NULL%1,  Shane%O’Mahony%NULL%1,  Carmen%Mikacenic%NULL%1]'
$ws.Cells.Item(4, 9).Value = '_PMC'

$ws.Cells.Item(5, 9).Value = ''

$ws.Cells.Item(6, 5).Value = '[NULL%NULL%NULL%0,  NULL%NULL%NULL%8,  NULL%NULL%NULL%0,  Stephanie%Bialek%NULL%1,  Ellen%Boundy%NULL%1,  Virginia%Bowen%NULL%1,  Nancy%Chow%NULL%1,  Amanda%Cohn%NULL%1,  Nicole%Dowling%NULL%1,  Sascha%Ellington%NULL%1,  Ryan%Gierke%NULL%1,  Aron%Hall%NULL%1,  Jessica%MacNeil%NULL%1,  Priti%Patel%NULL%1,  Georgina%Peacock%NULL%3,  Tamara%Pilishvili%NULL%1,  Hilda%Razzaghi%NULL%2,  Nia%Reed%NULL%1,  Matthew%Ritchey%NULL%1,  Erin%Sauber-Schatz%NULL%1]'
$ws.Cells.Item(6, 9).Value = '_PMC'

$ws.Cells.Item(7, 5).Value = '[Jianlei%Cao%NULL%0,  Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%1,  Wenlin%Cheng%NULL%1,  Lei%Yu%NULL%1,  Ya-Kun%Liu%NULL%1,  Xiaoyong%Hu%NULL%1,  Qiang%Liu%NULL%1]'
$ws.Cells.Item(7, 9).Value = '_PMC'

$ws.Cells.Item(8, 5).Value = '[Jun%Chen%NULL%0,  Tangkai%Qi%NULL%1,  Li%Liu%NULL%1,  Yun%Ling%NULL%1,  Zhiping%Qian%NULL%1,  Tao%Li%NULL%2,  Feng%Li%NULL%1,  Qingnian%Xu%NULL%1,  Yuyi%Zhang%NULL%1,  Shuibao%Xu%NULL%1,  Zhigang%Song%NULL%1,  Yigang%Zeng%NULL%1,  Yinzhong%Shen%NULL%1,  Yuxin%Shi%NULL%1,  Tongyu%Zhu%NULL%1,  Hongzhou%Lu%NULL%1]'
$ws.Cells.Item(8, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(9, 5).Value = '[Na%Zhu%NULL%0,  Dingyu%Zhang%NULL%2,  Wenling%Wang%NULL%1,  Xingwang%Li%NULL%1,  Bo%Yang%NULL%1,  Jingdong%Song%NULL%1,  Xiang%Zhao%NULL%1,  Baoying%Huang%NULL%1,  Weifeng%Shi%NULL%1,  Roujian%Lu%NULL%1,  Peihua%Niu%NULL%1,  Faxian%Zhan%NULL%1,  Xuejun%Ma%NULL%1,  Dayan%Wang%NULL%1,  Wenbo%Xu%NULL%1,  Guizhen%Wu%NULL%1,  George F.%Gao%NULL%1,  Wenjie%Tan%NULL%1]'
$ws.Cells.Item(9, 9).Value = '_PMC'

$ws.Cells.Item(10, 5).Value = '[Christopher M%Petrilli%NULL%5,  Simon A%Jones%NULL%1,  Jie%Yang%NULL%1,  Harish%Rajagopalan%NULL%1,  Luke%O’Donnell%NULL%1,  Yelena%Chernyak%NULL%1,  Katie A%Tobin%NULL%1,  Robert J%Cerfolio%NULL%1,  Fritz%Francois%NULL%2,  Leora I%Horwitz%NULL%1]'
$ws.Cells.Item(10, 9).Value = '_PMC'

$ws.Cells.Item(11, 9).Value = ''

$ws.Cells.Item(12, 5).Value = '[Arthur%Simonnet%NULL%0,  Mikael%Chetboun%NULL%1,  Julien%Poissy%NULL%1,  Violeta%Raverdy%NULL%1,  Jerome%Noulette%NULL%2,  Jerome%Noulette%NULL%0,  Alain%Duhamel%NULL%1,  Julien%Labreuche%NULL%1,  Daniel%Mathieu%NULL%2,  Francois%Pattou%francois.pattou@univ-lille.fr%1,  Merce%Jourdain%mercedes.jourdain@univ-lille.fr%2,  Merce%Jourdain%mercedes.jourdain@univ-lille.fr%0,  Robert%Caizzo%NULL%2,  Robert%Caizzo%NULL%0,  Morgan%Caplan%NULL%1,  Nicolas%Cousin%NULL%1,  Thibault%Duburcq%NULL%1,  Arthur%Durand%NULL%1,  Ahmed%El kalioubie%NULL%1,  Raphael%Favory%NULL%1,  Bruno%Garcia%NULL%1,  Patrick%Girardie%NULL%1,  Julien%Goutay%NULL%1,  Marion%Houard%NULL%1,  Emmanuelle%Jaillette%NULL%1,  Nicolas%Kostuj%NULL%1,  Geoffrey%Ledoux%NULL%1,  Daniel%Mathieu%NULL%0,  Anne Sophie%Moreau%NULL%1,  Christopher%Niles%NULL%1,  Saad%Nseir%NULL%1,  Thierry%Onimus%NULL%1,  Erika%Parmentier%NULL%1,  Sebastien%Préau%NULL%1,  Laurent%Robriquet%NULL%1,  Anahita%Rouze%NULL%1,  Sophie%Six%NULL%1,  Hélène%Verkindt%NULL%1]'
$ws.Cells.Item(12, 9).Value = '_PMC'

$ws.Cells.Item(13, 9).Value = ''

$ws.Cells.Item(14, 9).Value = ''

$ws.Cells.Item(15, 5).Value = '[Cheng-Yu%Chen%NULL%1,  Chen-Hsen%Lee%NULL%1,  Cheng-Yi%Liu%NULL%1,  Jia-Horng%Wang%NULL%1,  Lee-Min%Wang%NULL%1,  Reury-Perng%Perng%NULL%1]'
$ws.Cells.Item(15, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(16, 5).Value = '[Jiqian%Xu%NULL%1,  Xiaobo%Yang%NULL%1,  Luyu%Yang%NULL%1,  Xiaojing%Zou%NULL%1,  Yaxin%Wang%NULL%1,  Yongran%Wu%NULL%1,  Ting%Zhou%NULL%1,  Yin%Yuan%NULL%1,  Hong%Qi%NULL%1,  Shouzhi%Fu%NULL%1,  Hong%Liu%NULL%1,  Jia’an%Xia%NULL%1,  Zhengqin%Xu%NULL%1,  Yuan%Yu%NULL%1,  Ruiting%Li%NULL%0,  Yaqi%Ouyang%NULL%1,  Rui%Wang%NULL%1,  Lehao%Ren%NULL%1,  Yingying%Hu%NULL%1,  Dan%Xu%NULL%1,  Xin%Zhao%NULL%1,  Shiying%Yuan%shiying_yuan@163.com%1,  Dingyu%Zhang%1813886398@qq.com%0,  You%Shang%you_shanghust@163.com%1]'
$ws.Cells.Item(16, 9).Value = '_PMC_Springer'

$ws.Cells.Item(17, 9).Value = ''

$ws.Cells.Item(18, 5).Value = '[Wei-jie%Guan%NULL%0,  Wen-hua%Liang%NULL%1,  Yi%Zhao%NULL%1,  Heng-rui%Liang%NULL%1,  Zi-sheng%Chen%NULL%1,  Yi-min%Li%NULL%1,  Xiao-qing%Liu%NULL%1,  Ru-chong%Chen%NULL%1,  Chun-li%Tang%NULL%1,  Tao%Wang%NULL%1,  Chun-quan%Ou%NULL%1,  Li%Li%NULL%0,  Ping-yan%Chen%NULL%1,  Ling%Sang%NULL%1,  Wei%Wang%NULL%1,  Jian-fu%Li%NULL%1,  Cai-chen%Li%NULL%1,  Li-min%Ou%NULL%1,  Bo%Cheng%NULL%1,  Shan%Xiong%NULL%1,  Zheng-yi%Ni%NULL%1,  Jie%Xiang%NULL%2,  Yu%Hu%NULL%1,  Lei%Liu%NULL%1,  Hong%Shan%NULL%1,  Chun-liang%Lei%NULL%1,  Yi-xiang%Peng%NULL%1,  Li%Wei%NULL%1,  Yong%Liu%NULL%1,  Ya-hua%Hu%NULL%1,  Peng%Peng%NULL%2,  Jian-ming%Wang%NULL%1,  Ji-yang%Liu%NULL%1,  Zhong%Chen%NULL%1,  Gang%Li%NULL%1,  Zhi-jian%Zheng%NULL%1,  Shao-qin%Qiu%NULL%1,  Jie%Luo%NULL%1,  Chang-jiang%Ye%NULL%1,  Shao-yong%Zhu%NULL%1,  Lin-ling%Cheng%NULL%1,  Feng%Ye%NULL%1,  Shi-yue%Li%NULL%1,  Jin-ping%Zheng%NULL%1,  Nuo-fu%Zhang%NULL%1,  Nan-shan%Zhong%NULL%1,  Jian-xing%He%NULL%1]'
$ws.Cells.Item(18, 9).Value = '_PMC'

$ws.Cells.Item(19, 5).Value = '[Fei%Zhou%NULL%0,  Ting%Yu%NULL%2,  Ronghui%Du%NULL%1,  Guohui%Fan%NULL%1,  Ying%Liu%NULL%2,  Zhibo%Liu%NULL%1,  Jie%Xiang%NULL%0,  Yeming%Wang%NULL%1,  Bin%Song%NULL%1,  Xiaoying%Gu%NULL%1,  Lulu%Guan%NULL%1,  Yuan%Wei%NULL%2,  Hui%Li%NULL%1,  Xudong%Wu%NULL%1,  Jiuyang%Xu%NULL%1,  Shengjin%Tu%NULL%1,  Yi%Zhang%NULL%1,  Hua%Chen%NULL%1,  Bin%Cao%NULL%1]'
$ws.Cells.Item(19, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(20, 5).Value = '[Corrado%Lodigiani%NULL%1,  Giacomo%Iapichino%NULL%1,  Luca%Carenzo%NULL%1,  Maurizio%Cecconi%NULL%1,  Paola%Ferrazzi%NULL%1,  Tim%Sebastian%NULL%1,  Nils%Kucher%NULL%1,  Jan-Dirk%Studt%NULL%1,  Clara%Sacco%NULL%1,  Alexia%Bertuzzi%NULL%1,  Maria Teresa%Sandri%NULL%1,  Stefano%Barco%NULL%1,  NULL%NULL%NULL%0]'
$ws.Cells.Item(20, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(21, 5).Value = '[F.A.%Klok%NULL%1,  M.J.H.A.%Kruip%NULL%1,  N.J.M.%van der Meer%NULL%1,  M.S.%Arbous%NULL%1,  D.A.M.P.J.%Gommers%NULL%1,  K.M.%Kant%NULL%1,  F.H.J.%Kaptein%NULL%1,  J.%van Paassen%NULL%1,  M.A.M.%Stals%NULL%1,  M.V.%Huisman%NULL%1,  H.%Endeman%NULL%1]'
$ws.Cells.Item(21, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(22, 5).Value = '[Shaoqing%Lei%NULL%3,  Fang%Jiang%NULL%1,  Wating%Su%NULL%1,  Chang%Chen%NULL%1,  Jingli%Chen%NULL%1,  Wei%Mei%NULL%1,  Li-Ying%Zhan%NULL%1,  Yifan%Jia%NULL%1,  Liangqing%Zhang%NULL%1,  Danyong%Liu%NULL%1,  Zhong-Yuan%Xia%xiazhongyuan2005@aliyun.com%1,  Zhengyuan%Xia%zyxia@hkucc.hku.hk%1]'
$ws.Cells.Item(22, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(23, 5).Value = '[Annemarie B%Docherty%NULL%1,  Ewen M%Harrison%NULL%1,  Christopher A%Green%NULL%1,  Hayley E%Hardwick%NULL%1,  Riinu%Pius%NULL%1,  Lisa%Norman%NULL%1,  Karl A%Holden%NULL%1,  Jonathan M%Read%NULL%1,  Frank%Dondelinger%NULL%1,  Gail%Carson%NULL%1,  Laura%Merson%NULL%1,  James%Lee%NULL%1,  Daniel%Plotkin%NULL%1,  Louise%Sigfrid%NULL%1,  Sophie%Halpin%NULL%1,  Clare%Jackson%NULL%1,  Carrol%Gamble%NULL%1,  Peter W%Horby%NULL%1,  Jonathan S%Nguyen-Van-Tam%NULL%1,  Antonia%Ho%NULL%1,  Clark D%Russell%NULL%1,  Jake%Dunning%NULL%1,  Peter JM%Openshaw%NULL%1,  J Kenneth%Baillie%NULL%1,  Malcolm G%Semple%NULL%1]'
$ws.Cells.Item(23, 9).Value = '_PMC'

$ws.Cells.Item(24, 5).Value = '[Rong-Hui%Du%NULL%0,  Li-Min%Liu%NULL%1,  Wen%Yin%NULL%1,  Wen%Wang%NULL%1,  Lu-Lu%Guan%NULL%1,  Ming-Li%Yuan%NULL%1,  Yu-Lei%Li%NULL%1,  Yi%Hu%NULL%1,  Xu-Yan%Li%NULL%1,  Bing%Sun%NULL%1,  Peng%Peng%NULL%0,  Huan-Zhong%Shi%NULL%1]'
$ws.Cells.Item(24, 9).Value = '_PMC'

$ws.Cells.Item(25, 5).Value = '[' + "`n" + 'Ming%Hu%NULL%1,  Qiang%Zhou%NULL%1,  Ruiqiang%Zheng%NULL%1,  Xuyan%Li%NULL%1,  Jianmin%Ling%NULL%1,  Yumei%Chen%NULL%1,  Jing%Jia%NULL%1,  Cuihong%Xie%xiecuihong08@163.com%1]'
$ws.Cells.Item(25, 9).Value = '_PMC_Springer'

$ws.Cells.Item(26, 5).Value = '[Evgeny V.%Fominskiy%NULL%1,  Anna Mara%Scandroglio%NULL%1,  Giacomo%Monti%NULL%1,  Maria Grazia%Calabrò%NULL%1,  Giovanni%Landoni%NULL%1,  Antonio%Dell''Acqua%NULL%1,  Luigi%Beretta%NULL%1,  Elena%Moizo%NULL%1,  Alfredo%Ravizza%NULL%1,  Fabrizio%Monaco%NULL%1,  Corrado%Campochiaro%NULL%1,  Marina%Pieri%NULL%1,  Maria Luisa%Azzolini%NULL%1,  Giovanni%Borghi%NULL%1,  Martina%Crivellari%NULL%1,  Caterina%Conte%NULL%1,  Cristina%Mattioli%NULL%1,  Paolo%Silvani%NULL%1,  Milena%Mucci%NULL%1,  Stefano%Turi%NULL%1,  Stefano%Tentori%NULL%1,  Martina%Baiardo Redaelli%NULL%1,  Marianna%Sartorelli%NULL%1,  Piera%Angelillo%NULL%1,  Alessandro%Belletti%NULL%1,  Pasquale%Nardelli%NULL%1,  Francesco Giuseppe%Nisi%NULL%1,  Gabriele%Valsecchi%NULL%1,  Cristina%Barberio%NULL%1,  Fabio%Ciceri%NULL%1,  Ary%Serpa Neto%NULL%1,  Lorenzo%Dagna%NULL%1,  Rinaldo%Bellomo%NULL%1,  Alberto%Zangrillo%NULL%1]'
$ws.Cells.Item(26, 9).Value = '_PMC'

$ws.Cells.Item(27, 9).Value = ''

$ws.Cells.Item(28, 5).Value = '[Nanshan%Chen%NULL%0,  Min%Zhou%NULL%1,  Xuan%Dong%NULL%1,  Jieming%Qu%NULL%1,  Fengyun%Gong%NULL%1,  Yang%Han%NULL%1,  Yang%Qiu%NULL%1,  Jingli%Wang%NULL%1,  Ying%Liu%NULL%0,  Yuan%Wei%NULL%0,  Jia''an%Xia%NULL%1,  Ting%Yu%NULL%0,  Xinxin%Zhang%NULL%1,  Li%Zhang%NULL%1]'
$ws.Cells.Item(28, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(29, 9).Value = ''

$ws.Cells.Item(30, 9).Value = ''

$ws.Cells.Item(31, 9).Value = ''

$ws.Cells.Item(32, 9).Value = ''

$ws.Cells.Item(33, 5).Value = '[Mustafa%Saad%NULL%1,  Ali S.%Omrani%NULL%1,  Kamran%Baig%NULL%1,  Abdelkarim%Bahloul%NULL%1,  Fatehi%Elzein%NULL%1,  Mohammad Abdul%Matin%NULL%1,  Mohei A.A.%Selim%NULL%1,  Mohammed Al%Mutairi%NULL%1,  Daifullah Al%Nakhli%NULL%1,  Amal Y. Al%Aidaroos%NULL%1,  Nisreen Al%Sherbeeni%NULL%1,  Hesham I.%Al-Khashan%NULL%1,  Ziad A.%Memish%zmemish@yahoo.com%1,  Ali M.%Albarrak%NULL%1]'
$ws.Cells.Item(33, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(34, 5).Value = '[Ashraf Abdel%Halim%NULL%1,  Badr%Alsayed%NULL%1,  Sameh%Embarak%NULL%1,  Taha%Yaseen%NULL%1,  Salwa%Dabbous%NULL%1]'
$ws.Cells.Item(34, 9).Value = '_PMC_elsevier'

$ws.Cells.Item(35, 5).Value = '[Musa A.%Garbati%NULL%1,  Shamsudeen F.%Fagbo%NULL%1,  Vicky J.%Fang%NULL%1,  Leila%Skakni%NULL%1,  Mercy%Joseph%NULL%1,  Tariq A.%Wani%NULL%1,  Benjamin J.%Cowling%NULL%1,  Malik%Peiris%NULL%1,  Ahmed%Hakawi%NULL%1,  Renee W.Y.%Chan%NULL%2,  Renee W.Y.%Chan%NULL%0]'
$ws.Cells.Item(35, 9).Value = '_PMC'

$ws.Cells.Item(36, 5).Value = '[Mohammed%Al Ghamdi%dr_mohd1391@yahoo.com%1,  Khalid M.%Alghamdi%K.m.a02@hotmail.com%1,  Yasmeen%Ghandoora%Yass.m.gh@hotmail.com%1,  Ameera%Alzahrani%Ameera.alzahrani@outlook.com%1,  Fatmah%Salah%Fs4l_med2011@yahoo.com%1,  Abdulmoatani%Alsulami%moatani@hotmail.com%1,  Mayada F.%Bawayan%mayada-f-b@hotmail.com%1,  Dhananjay%Vaidya%dvaidya@jhmi.edu%1,  Trish M.%Perl%tperl@jhmi.edu%1,  Geeta%Sood%gsood1@jhmi.edu%1]'
$ws.Cells.Item(36, 9).Value = '_PMC_Springer'

